$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G9").Value = 0.125390957193042
$ws.Range("G10").Value = 0.125390957193042
$ws.Range("G11").Value = 1.40652126935725
$ws.Range("G12").Value = 1.40652126935725
$ws.Range("F13").Value = 0.6838
$ws.Range("G13").Value = 1.51341694915254
$ws.Range("H13").Value = 9.5381
$ws.Range("F14").Value = 0.6838
$ws.Range("G14").Value = 1.51341694915254
$ws.Range("H14").Value = 9.5381
$ws.Range("G25").Value = 0.139613786341024
$ws.Range("G26").Value = 0.139613786341024
$ws.Range("G27").Value = 1.0507007152289
$ws.Range("G28").Value = 1.0507007152289
$ws.Range("F29").Value = 0.7146
$ws.Range("G29").Value = 1.3282593220339
$ws.Range("H29").Value = 9.5381
$ws.Range("L29").Value = 0.2362
$ws.Range("F30").Value = 0.7146
$ws.Range("G30").Value = 1.3282593220339
$ws.Range("H30").Value = 9.5381
$ws.Range("L30").Value = 0.2362
$ws.Range("G41").Value = 0.172323567152406
$ws.Range("G42").Value = 0.172323567152406
$ws.Range("G45").Value = 1.51814333333333
$ws.Range("H45").Value = 9.5381
$ws.Range("L45").Value = 0.34535
$ws.Range("N45").Value = 4.49959
$ws.Range("G46").Value = 1.51814333333333
$ws.Range("H46").Value = 9.5381
$ws.Range("L46").Value = 0.34535
$ws.Range("N46").Value = 4.49959
$ws.Range("G58").Value = 0.29119135698524
$ws.Range("G59").Value = 0.29119135698524
$ws.Range("G62").Value = 1.56992666666667
$ws.Range("H62").Value = 9.5381
$ws.Range("I62").Value = 4.38215
$ws.Range("L62").Value = 0.70745
$ws.Range("M62").Value = 2.73666
$ws.Range("G63").Value = 1.56992666666667
$ws.Range("H63").Value = 9.5381
$ws.Range("I63").Value = 4.38215
$ws.Range("L63").Value = 0.70745
$ws.Range("M63").Value = 2.73666
$ws.Range("G75").Value = 0.245082111675764
$ws.Range("G76").Value = 0.245082111675764
$ws.Range("G79").Value = 1.47502666666667
$ws.Range("H79").Value = 9.5381
$ws.Range("I79").Value = 3.91215
$ws.Range("L79").Value = 0.36965
$ws.Range("M79").Value = 2.73666
$ws.Range("G80").Value = 1.47502666666667
$ws.Range("H80").Value = 9.5381
$ws.Range("I80").Value = 3.91215
$ws.Range("L80").Value = 0.36965
$ws.Range("M80").Value = 2.73666
$ws.Range("G92").Value = 0.300494332972543
$ws.Range("G93").Value = 0.300494332972543
$ws.Range("G96").Value = 1.51539666666667
$ws.Range("L96").Value = 1.04455
$ws.Range("M96").Value = 2.73666
$ws.Range("G97").Value = 1.51539666666667
$ws.Range("L97").Value = 1.04455
$ws.Range("M97").Value = 2.73666
$ws.Range("G109").Value = 0.326280047014922
$ws.Range("G110").Value = 0.326280047014922
$ws.Range("F113").Value = 1.3661
$ws.Range("G113").Value = 1.605585
$ws.Range("I113").Value = 4.00715
$ws.Range("L113").Value = 1.04455
$ws.Range("M113").Value = 2.91389
$ws.Range("F114").Value = 1.3661
$ws.Range("G114").Value = 1.605585
$ws.Range("I114").Value = 4.00715
$ws.Range("L114").Value = 1.04455
$ws.Range("M114").Value = 2.91389
$ws.Range("G144").Value = 0.746616822399917
$ws.Range("G145").Value = 0.746616822399917
$ws.Range("G160").Value = 0.698890801740266
$ws.Range("G161").Value = 0.698890801740266
$ws.Cells.Item(168, 1).Value = "Rangitawa Stream at ds Halcombe oxpond"
$ws.Cells.Item(168, 2).Value = "DRP (95th Percentile)"
$ws.Cells.Item(168, 3).Value = "D"
$ws.Cells.Item(168, 4).Value = "2019 - 2023"
$ws.Cells.Item(168, 5).Value = "Impact"
$ws.Cells.Item(168, 6).Value = 0.045
$ws.Cells.Item(168, 7).Value = 0.127614035087719
$ws.Cells.Item(168, 8).Value = 1.57
$ws.Cells.Item(168, 9).Value = 0.6721
$ws.Cells.Item(168, 12).Value = 0.0555
$ws.Cells.Item(168, 13).Value = 0.12053
$ws.Cells.Item(168, 14).Value = 0.3725
$ws.Cells.Item(168, 15).Value = 1811681.523
$ws.Cells.Item(168, 16).Value = 5554500.976
$ws.Cells.Item(168, 17).Value = "Manawatu District"
$ws.Cells.Item(168, 18).Value = "Rangitīkei-Turakina"
$ws.Cells.Item(168, 19).Value = "Coastal Rangitikei"
$ws.Cells.Item(168, 20).Value = "Rang_4a"
$ws.Cells.Item(168, 21).Value = "mg/L"
$ws.Cells.Item(169, 1).Value = "Rangitawa Stream at ds Halcombe oxpond"
$ws.Cells.Item(169, 2).Value = "DRP (Median)"
$ws.Cells.Item(169, 3).Value = "D"
$ws.Cells.Item(169, 4).Value = "2019 - 2023"
$ws.Cells.Item(169, 5).Value = "Impact"
$ws.Cells.Item(169, 6).Value = 0.045
$ws.Cells.Item(169, 7).Value = 0.127614035087719
$ws.Cells.Item(169, 8).Value = 1.57
$ws.Cells.Item(169, 9).Value = 0.6721
$ws.Cells.Item(169, 12).Value = 0.0555
$ws.Cells.Item(169, 13).Value = 0.12053
$ws.Cells.Item(169, 14).Value = 0.3725
$ws.Cells.Item(169, 15).Value = 1811681.523
$ws.Cells.Item(169, 16).Value = 5554500.976
$ws.Cells.Item(169, 17).Value = "Manawatu District"
$ws.Cells.Item(169, 18).Value = "Rangitīkei-Turakina"
$ws.Cells.Item(169, 19).Value = "Coastal Rangitikei"
$ws.Cells.Item(169, 20).Value = "Rang_4a"
$ws.Cells.Item(169, 21).Value = "mg/L"
$ws.Cells.Item(170, 1).Value = "Rangitawa Stream at ds Halcombe oxpond"
$ws.Cells.Item(170, 2).Value = "E coli (>260)"
$ws.Cells.Item(170, 3).Value = "E"
$ws.Cells.Item(170, 4).Value = "2019 - 2023"
$ws.Cells.Item(170, 5).Value = "Impact"
$ws.Cells.Item(170, 6).Value = 550
$ws.Cells.Item(170, 7).Value = 2143.84210526316
$ws.Cells.Item(170, 8).Value = 41000
$ws.Cells.Item(170, 9).Value = 7561.2
$ws.Cells.Item(170, 10).Value = 50.8771929824561
$ws.Cells.Item(170, 11).Value = 71.92982456140351
$ws.Cells.Item(170, 12).Value = 1115
$ws.Cells.Item(170, 13).Value = 2424
$ws.Cells.Item(170, 14).Value = 4800
$ws.Cells.Item(170, 15).Value = 1811681.523
$ws.Cells.Item(170, 16).Value = 5554500.976
$ws.Cells.Item(170, 17).Value = "Manawatu District"
$ws.Cells.Item(170, 18).Value = "Rangitīkei-Turakina"
$ws.Cells.Item(170, 19).Value = "Coastal Rangitikei"
$ws.Cells.Item(170, 20).Value = "Rang_4a"
$ws.Cells.Item(170, 21).Value = "% exceedances over 260/100 mL"
$ws.Cells.Item(171, 1).Value = "Rangitawa Stream at ds Halcombe oxpond"
$ws.Cells.Item(171, 2).Value = "E coli (>540)"
$ws.Cells.Item(171, 3).Value = "E"
$ws.Cells.Item(171, 4).Value = "2019 - 2023"
$ws.Cells.Item(171, 5).Value = "Impact"
$ws.Cells.Item(171, 6).Value = 550
$ws.Cells.Item(171, 7).Value = 2143.84210526316
$ws.Cells.Item(171, 8).Value = 41000
$ws.Cells.Item(171, 9).Value = 7561.2
$ws.Cells.Item(171, 10).Value = 50.8771929824561
$ws.Cells.Item(171, 11).Value = 71.92982456140351
$ws.Cells.Item(171, 12).Value = 1115
$ws.Cells.Item(171, 13).Value = 2424
$ws.Cells.Item(171, 14).Value = 4800
$ws.Cells.Item(171, 15).Value = 1811681.523
$ws.Cells.Item(171, 16).Value = 5554500.976
$ws.Cells.Item(171, 17).Value = "Manawatu District"
$ws.Cells.Item(171, 18).Value = "Rangitīkei-Turakina"
$ws.Cells.Item(171, 19).Value = "Coastal Rangitikei"
$ws.Cells.Item(171, 20).Value = "Rang_4a"
$ws.Cells.Item(171, 21).Value = "% exceedances over 540/100 mL"
$ws.Cells.Item(172, 1).Value = "Rangitawa Stream at ds Halcombe oxpond"
$ws.Cells.Item(172, 2).Value = "E coli (Median)"
$ws.Cells.Item(172, 3).Value = "E"
$ws.Cells.Item(172, 4).Value = "2019 - 2023"
$ws.Cells.Item(172, 5).Value = "Impact"
$ws.Cells.Item(172, 6).Value = 550
$ws.Cells.Item(172, 7).Value = 2143.84210526316
$ws.Cells.Item(172, 8).Value = 41000
$ws.Cells.Item(172, 9).Value = 7561.2
$ws.Cells.Item(172, 10).Value = 50.8771929824561
$ws.Cells.Item(172, 11).Value = 71.92982456140351
$ws.Cells.Item(172, 12).Value = 1115
$ws.Cells.Item(172, 13).Value = 2424
$ws.Cells.Item(172, 14).Value = 4800
$ws.Cells.Item(172, 15).Value = 1811681.523
$ws.Cells.Item(172, 16).Value = 5554500.976
$ws.Cells.Item(172, 17).Value = "Manawatu District"
$ws.Cells.Item(172, 18).Value = "Rangitīkei-Turakina"
$ws.Cells.Item(172, 19).Value = "Coastal Rangitikei"
$ws.Cells.Item(172, 20).Value = "Rang_4a"
$ws.Cells.Item(172, 21).Value = "E. coli/100 mL"
$ws.Cells.Item(173, 1).Value = "Rangitawa Stream at ds Halcombe oxpond"
$ws.Cells.Item(173, 2).Value = "E coli (95th Percentile)"
$ws.Cells.Item(173, 3).Value = "E"
$ws.Cells.Item(173, 4).Value = "2019 - 2023"
$ws.Cells.Item(173, 5).Value = "Impact"
$ws.Cells.Item(173, 6).Value = 550
$ws.Cells.Item(173, 7).Value = 2143.84210526316
$ws.Cells.Item(173, 8).Value = 41000
$ws.Cells.Item(173, 9).Value = 7561.2
$ws.Cells.Item(173, 10).Value = 50.8771929824561
$ws.Cells.Item(173, 11).Value = 71.92982456140351
$ws.Cells.Item(173, 12).Value = 1115
$ws.Cells.Item(173, 13).Value = 2424
$ws.Cells.Item(173, 14).Value = 4800
$ws.Cells.Item(173, 15).Value = 1811681.523
$ws.Cells.Item(173, 16).Value = 5554500.976
$ws.Cells.Item(173, 17).Value = "Manawatu District"
$ws.Cells.Item(173, 18).Value = "Rangitīkei-Turakina"
$ws.Cells.Item(173, 19).Value = "Coastal Rangitikei"
$ws.Cells.Item(173, 20).Value = "Rang_4a"
$ws.Cells.Item(173, 21).Value = "E. coli/100 mL"
$ws.Cells.Item(174, 1).Value = "Rangitawa Stream at ds Halcombe oxpond"
$ws.Cells.Item(174, 2).Value = "Ammoniacal-N (95th Percentile)"
$ws.Cells.Item(174, 3).Value = "C"
$ws.Cells.Item(174, 4).Value = "2019 - 2023"
$ws.Cells.Item(174, 5).Value = "Impact"
$ws.Cells.Item(174, 6).Value = 0.06222
$ws.Cells.Item(174, 7).Value = 0.249785303421946
$ws.Cells.Item(174, 8).Value = 2.71906648967442
$ws.Cells.Item(174, 9).Value = 1.03894
$ws.Cells.Item(174, 12).Value = 0.09723999999999999
$ws.Cells.Item(174, 13).Value = 0.55171
$ws.Cells.Item(174, 14).Value = 0.9129699999999999
$ws.Cells.Item(174, 15).Value = 1811681.523
$ws.Cells.Item(174, 16).Value = 5554500.976
$ws.Cells.Item(174, 17).Value = "Manawatu District"
$ws.Cells.Item(174, 18).Value = "Rangitīkei-Turakina"
$ws.Cells.Item(174, 19).Value = "Coastal Rangitikei"
$ws.Cells.Item(174, 20).Value = "Rang_4a"
$ws.Cells.Item(174, 21).Value = "mg NH4-N/L"
$ws.Cells.Item(175, 1).Value = "Rangitawa Stream at ds Halcombe oxpond"
$ws.Cells.Item(175, 2).Value = "Ammoniacal-N (Median)"
$ws.Cells.Item(175, 3).Value = "B"
$ws.Cells.Item(175, 4).Value = "2019 - 2023"
$ws.Cells.Item(175, 5).Value = "Impact"
$ws.Cells.Item(175, 6).Value = 0.06222
$ws.Cells.Item(175, 7).Value = 0.249785303421946
$ws.Cells.Item(175, 8).Value = 2.71906648967442
$ws.Cells.Item(175, 9).Value = 1.03894
$ws.Cells.Item(175, 12).Value = 0.09723999999999999
$ws.Cells.Item(175, 13).Value = 0.55171
$ws.Cells.Item(175, 14).Value = 0.9129699999999999
$ws.Cells.Item(175, 15).Value = 1811681.523
$ws.Cells.Item(175, 16).Value = 5554500.976
$ws.Cells.Item(175, 17).Value = "Manawatu District"
$ws.Cells.Item(175, 18).Value = "Rangitīkei-Turakina"
$ws.Cells.Item(175, 19).Value = "Coastal Rangitikei"
$ws.Cells.Item(175, 20).Value = "Rang_4a"
$ws.Cells.Item(175, 21).Value = "mg NH4-N/L"
$ws.Cells.Item(176, 1).Value = "Rangitawa Stream at ds Halcombe oxpond"
$ws.Cells.Item(176, 2).Value = "Nitrate-N (95th Percentile)"
$ws.Cells.Item(176, 3).Value = "B"
$ws.Cells.Item(176, 4).Value = "2019 - 2023"
$ws.Cells.Item(176, 5).Value = "Impact"
$ws.Cells.Item(176, 6).Value = 0.251
$ws.Cells.Item(176, 7).Value = 0.674732907003424
$ws.Cells.Item(176, 8).Value = 4.64
$ws.Cells.Item(176, 9).Value = 3.1885
$ws.Cells.Item(176, 12).Value = 0.172
$ws.Cells.Item(176, 13).Value = 1.1305
$ws.Cells.Item(176, 14).Value = 2.531
$ws.Cells.Item(176, 15).Value = 1811681.523
$ws.Cells.Item(176, 16).Value = 5554500.976
$ws.Cells.Item(176, 17).Value = "Manawatu District"
$ws.Cells.Item(176, 18).Value = "Rangitīkei-Turakina"
$ws.Cells.Item(176, 19).Value = "Coastal Rangitikei"
$ws.Cells.Item(176, 20).Value = "Rang_4a"
$ws.Cells.Item(176, 21).Value = "mg NO3-N/L"
$ws.Cells.Item(177, 1).Value = "Rangitawa Stream at ds Halcombe oxpond"
$ws.Cells.Item(177, 2).Value = "Nitrate-N (Median)"
$ws.Cells.Item(177, 3).Value = "A"
$ws.Cells.Item(177, 4).Value = "2019 - 2023"
$ws.Cells.Item(177, 5).Value = "Impact"
$ws.Cells.Item(177, 6).Value = 0.251
$ws.Cells.Item(177, 7).Value = 0.674732907003424
$ws.Cells.Item(177, 8).Value = 4.64
$ws.Cells.Item(177, 9).Value = 3.1885
$ws.Cells.Item(177, 12).Value = 0.172
$ws.Cells.Item(177, 13).Value = 1.1305
$ws.Cells.Item(177, 14).Value = 2.531
$ws.Cells.Item(177, 15).Value = 1811681.523
$ws.Cells.Item(177, 16).Value = 5554500.976
$ws.Cells.Item(177, 17).Value = "Manawatu District"
$ws.Cells.Item(177, 18).Value = "Rangitīkei-Turakina"
$ws.Cells.Item(177, 19).Value = "Coastal Rangitikei"
$ws.Cells.Item(177, 20).Value = "Rang_4a"
$ws.Cells.Item(177, 21).Value = "mg NO3-N/L"
$ws.Cells.Item(178, 1).Value = "Rangitawa Stream at ds Halcombe oxpond"
$ws.Cells.Item(178, 2).Value = "Soluble Inorganic Nitrogen (95th Percentile)"
$ws.Cells.Item(178, 4).Value = "2019 - 2023"
$ws.Cells.Item(178, 5).Value = "Impact"
$ws.Cells.Item(178, 6).Value = 0.487
$ws.Cells.Item(178, 7).Value = 1.26859649122807
$ws.Cells.Item(178, 8).Value = 7.11
$ws.Cells.Item(178, 9).Value = 4.032
$ws.Cells.Item(178, 12).Value = 0.4645
$ws.Cells.Item(178, 13).Value = 2.6362
$ws.Cells.Item(178, 14).Value = 3.5187
$ws.Cells.Item(178, 15).Value = 1811681.523
$ws.Cells.Item(178, 16).Value = 5554500.976
$ws.Cells.Item(178, 17).Value = "Manawatu District"
$ws.Cells.Item(178, 18).Value = "Rangitīkei-Turakina"
$ws.Cells.Item(178, 19).Value = "Coastal Rangitikei"
$ws.Cells.Item(178, 20).Value = "Rang_4a"
$ws.Cells.Item(178, 21).Value = "g/m3"
$ws.Cells.Item(179, 1).Value = "Rangitawa Stream at ds Halcombe oxpond"
$ws.Cells.Item(179, 2).Value = "Soluble Inorganic Nitrogen (Median)"
$ws.Cells.Item(179, 4).Value = "2019 - 2023"
$ws.Cells.Item(179, 5).Value = "Impact"
$ws.Cells.Item(179, 6).Value = 0.487
$ws.Cells.Item(179, 7).Value = 1.26859649122807
$ws.Cells.Item(179, 8).Value = 7.11
$ws.Cells.Item(179, 9).Value = 4.032
$ws.Cells.Item(179, 12).Value = 0.4645
$ws.Cells.Item(179, 13).Value = 2.6362
$ws.Cells.Item(179, 14).Value = 3.5187
$ws.Cells.Item(179, 15).Value = 1811681.523
$ws.Cells.Item(179, 16).Value = 5554500.976
$ws.Cells.Item(179, 17).Value = "Manawatu District"
$ws.Cells.Item(179, 18).Value = "Rangitīkei-Turakina"
$ws.Cells.Item(179, 19).Value = "Coastal Rangitikei"
$ws.Cells.Item(179, 20).Value = "Rang_4a"
$ws.Cells.Item(179, 21).Value = "g/m3"
$ws.Cells.Item(180, 1).Value = "Rangitawa Stream at ds Halcombe oxpond"
$ws.Cells.Item(180, 2).Value = "Total Nitrogen (95th Percentile)"
$ws.Cells.Item(180, 4).Value = "2019 - 2023"
$ws.Cells.Item(180, 5).Value = "Impact"
$ws.Cells.Item(180, 6).Value = 1.26
$ws.Cells.Item(180, 7).Value = 2.03333333333333
$ws.Cells.Item(180, 8).Value = 8.619999999999999
$ws.Cells.Item(180, 9).Value = 5.218
$ws.Cells.Item(180, 12).Value = 1.105
$ws.Cells.Item(180, 13).Value = 3.5407
$ws.Cells.Item(180, 14).Value = 4.7076
$ws.Cells.Item(180, 15).Value = 1811681.523
$ws.Cells.Item(180, 16).Value = 5554500.976
$ws.Cells.Item(180, 17).Value = "Manawatu District"
$ws.Cells.Item(180, 18).Value = "Rangitīkei-Turakina"
$ws.Cells.Item(180, 19).Value = "Coastal Rangitikei"
$ws.Cells.Item(180, 20).Value = "Rang_4a"
$ws.Cells.Item(180, 21).Value = "g/m3"
$ws.Cells.Item(181, 1).Value = "Rangitawa Stream at ds Halcombe oxpond"
$ws.Cells.Item(181, 2).Value = "Total Nitrogen (Median)"
$ws.Cells.Item(181, 4).Value = "2019 - 2023"
$ws.Cells.Item(181, 5).Value = "Impact"
$ws.Cells.Item(181, 6).Value = 1.26
$ws.Cells.Item(181, 7).Value = 2.03333333333333
$ws.Cells.Item(181, 8).Value = 8.619999999999999
$ws.Cells.Item(181, 9).Value = 5.218
$ws.Cells.Item(181, 12).Value = 1.105
$ws.Cells.Item(181, 13).Value = 3.5407
$ws.Cells.Item(181, 14).Value = 4.7076
$ws.Cells.Item(181, 15).Value = 1811681.523
$ws.Cells.Item(181, 16).Value = 5554500.976
$ws.Cells.Item(181, 17).Value = "Manawatu District"
$ws.Cells.Item(181, 18).Value = "Rangitīkei-Turakina"
$ws.Cells.Item(181, 19).Value = "Coastal Rangitikei"
$ws.Cells.Item(181, 20).Value = "Rang_4a"
$ws.Cells.Item(181, 21).Value = "g/m3"
$ws.Cells.Item(182, 1).Value = "Rangitawa Stream at ds Halcombe oxpond"
$ws.Cells.Item(182, 2).Value = "Total Phosphorus (95th Percentile)"
$ws.Cells.Item(182, 4).Value = "2019 - 2023"
$ws.Cells.Item(182, 5).Value = "Impact"
$ws.Cells.Item(182, 6).Value = 0.139
$ws.Cells.Item(182, 7).Value = 0.254842105263158
$ws.Cells.Item(182, 8).Value = 1.97
$ws.Cells.Item(182, 9).Value = 0.98445
$ws.Cells.Item(182, 12).Value = 0.159
$ws.Cells.Item(182, 13).Value = 0.37701
$ws.Cells.Item(182, 14).Value = 0.5508999999999999
$ws.Cells.Item(182, 15).Value = 1811681.523
$ws.Cells.Item(182, 16).Value = 5554500.976
$ws.Cells.Item(182, 17).Value = "Manawatu District"
$ws.Cells.Item(182, 18).Value = "Rangitīkei-Turakina"
$ws.Cells.Item(182, 19).Value = "Coastal Rangitikei"
$ws.Cells.Item(182, 20).Value = "Rang_4a"
$ws.Cells.Item(182, 21).Value = "g/m3"
$ws.Cells.Item(183, 1).Value = "Rangitawa Stream at ds Halcombe oxpond"
$ws.Cells.Item(183, 2).Value = "Total Phosphorus (Median)"
$ws.Cells.Item(183, 4).Value = "2019 - 2023"
$ws.Cells.Item(183, 5).Value = "Impact"
$ws.Cells.Item(183, 6).Value = 0.139
$ws.Cells.Item(183, 7).Value = 0.254842105263158
$ws.Cells.Item(183, 8).Value = 1.97
$ws.Cells.Item(183, 9).Value = 0.98445
$ws.Cells.Item(183, 12).Value = 0.159
$ws.Cells.Item(183, 13).Value = 0.37701
$ws.Cells.Item(183, 14).Value = 0.5508999999999999
$ws.Cells.Item(183, 15).Value = 1811681.523
$ws.Cells.Item(183, 16).Value = 5554500.976
$ws.Cells.Item(183, 17).Value = "Manawatu District"
$ws.Cells.Item(183, 18).Value = "Rangitīkei-Turakina"
$ws.Cells.Item(183, 19).Value = "Coastal Rangitikei"
$ws.Cells.Item(183, 20).Value = "Rang_4a"
$ws.Cells.Item(183, 21).Value = "g/m3"
